# V 0.48-B45 PreRel for Audio Init Tests
# - Add AP Attitude Hold item ("ATT") as a new column on the Tabelle2
#   export/merge sheet, right before the existing END_OF_COL marker.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

# Insert a brand new column at DQ - this shifts the old DQ (END_OF_COL
# marker, value "X") to DR and the old DR (Title formula column) to DS,
# exactly like Excel's normal "Insert Column" behavior (styles/widths
# move with the shifted cells).
$ws.Columns("DQ").Insert()

# Header for the newly inserted column.
$ws.Range("DQ1").Value = "ATT"

# Body rows use the same "|" placeholder used by all the other
# ANI/flag columns in this sheet (DC:DP).
$ws.Range("DQ2:DQ40").Value = "|"

# Reflect the cursor position left behind by the edit.
$ws.Range("DO42").Select()
